# Update the "想去人数" (F column) figures on the 展览 and 全部类型 sheets
# to match the newly generated gh-pages data snapshot.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 157
$ws1.Range("F5").Value  = 1791
$ws1.Range("F7").Value  = 15
$ws1.Range("F8").Value  = 156
$ws1.Range("F9").Value  = 2143
$ws1.Range("F10").Value = 44
$ws1.Range("F11").Value = 148
$ws1.Range("F12").Value = 1358
$ws1.Range("F23").Value = 1161
$ws1.Range("F24").Value = 8
$ws1.Range("F25").Value = 345
$ws1.Range("F27").Value = 274
$ws1.Range("F28").Value = 341

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 157
$ws4.Range("F5").Value  = 1791
$ws4.Range("F8").Value  = 15
$ws4.Range("F9").Value  = 156
$ws4.Range("F10").Value = 2143
$ws4.Range("F11").Value = 44
$ws4.Range("F12").Value = 148
$ws4.Range("F13").Value = 1358
$ws4.Range("F24").Value = 1161
$ws4.Range("F25").Value = 8
$ws4.Range("F26").Value = 345
$ws4.Range("F28").Value = 274
$ws4.Range("F29").Value = 341

$wb.Save()
